$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary row 12: average of column J (J2:J11)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary rows 14-17: labels in column A, aggregate formulas in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Formatting: bold, larger, vertically centered labels for the B14:B17 summary values
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B15:B17").PasteSpecial(-4122) | Out-Null

# Formatting: bold for the J12 summary average
$ws.Range("J12").Font.Bold = $true

# Page setup: A4 portrait (as produced by the resave)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection left on J12, matching the saved UI state
$ws.Range("J12").Select() | Out-Null
